$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.5725316666666668
$ws.Range("H2").Value = 1.717595
$ws.Range("I2").Value = 0.3864899584549088
$ws.Range("J2").Value = 0.3864899584549088
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 6.322177333333333
$ws.Range("N2").Value = 18.966532
$ws.Range("O2").Value = 0.08271011762055308
$ws.Range("P2").Value = 0.0827101176205531
$ws.Range("Q2").Value = 3.619646725615556
$ws.Range("R2").Value = 32.57682053054
$ws.Range("S2").Value = 0.03196662992296818
$ws.Range("T2").Value = 0.03196662992296818

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.5725316666666668
$ws.Range("H3").Value = 1.717595
$ws.Range("I3").Value = 0.3864899584549088
$ws.Range("J3").Value = 0.3864899584549088
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 41.286995
$ws.Range("N3").Value = 123.860985
$ws.Range("O3").Value = 0.5401386314560596
$ws.Range("P3").Value = 0.5401386314560597
$ws.Range("Q3").Value = 23.63811205900834
$ws.Range("R3").Value = 212.743008531075
$ws.Range("S3").Value = 0.2087581572313437
$ws.Range("T3").Value = 0.2087581572313438

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.5725316666666668
$ws.Range("H4").Value = 1.717595
$ws.Range("I4").Value = 0.3864899584549088
$ws.Range("J4").Value = 0.3864899584549088
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 27.73243066666667
$ws.Range("N4").Value = 83.197292
$ws.Range("O4").Value = 0.3628105447549136
$ws.Range("P4").Value = 0.3628105447549136
$ws.Range("Q4").Value = 15.87769475030445
$ws.Range("R4").Value = 142.89925275274
$ws.Range("S4").Value = 0.1402226323693294
$ws.Range("T4").Value = 0.1402226323693294

$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.5725316666666668
$ws.Range("H5").Value = 1.717595
$ws.Range("I5").Value = 0.3864899584549088
$ws.Range("J5").Value = 0.3864899584549088
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.096171666666667
$ws.Range("N5").Value = 3.288515
$ws.Range("O5").Value = 0.01434070616847367
$ws.Range("P5").Value = 0.01434070616847367
$ws.Range("Q5").Value = 0.6275929912694446
$ws.Range("R5").Value = 5.648336921425001
$ws.Range("S5").Value = 0.005542538931267441
$ws.Range("T5").Value = 0.005542538931267441

$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.3064436666666667
$ws.Range("H6").Value = 0.919331
$ws.Range("I6").Value = 0.2068661122070742
$ws.Range("J6").Value = 0.2068661122070743
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 6.322177333333333
$ws.Range("N6").Value = 18.966532
$ws.Range("O6").Value = 0.08271011762055308
$ws.Range("P6").Value = 0.0827101176205531
$ws.Range("Q6").Value = 1.937391203343556
$ws.Range("R6").Value = 17.436520830092
$ws.Range("S6").Value = 0.01710992047235364
$ws.Range("T6").Value = 0.01710992047235365

$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.3064436666666667
$ws.Range("H7").Value = 0.919331
$ws.Range("I7").Value = 0.2068661122070742
$ws.Range("J7").Value = 0.2068661122070743
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 41.286995
$ws.Range("N7").Value = 123.860985
$ws.Range("O7").Value = 0.5401386314560596
$ws.Range("P7").Value = 0.5401386314560597
$ws.Range("Q7").Value = 12.65213813344833
$ws.Range("R7").Value = 113.869243201035
$ws.Range("S7").Value = 0.1117363787421647
$ws.Range("T7").Value = 0.1117363787421648

$ws.Range("E8").Value = 2
$ws.Range("F8").Value = 0.6666666666666666
$ws.Range("G8").Value = 0.3064436666666667
$ws.Range("H8").Value = 0.919331
$ws.Range("I8").Value = 0.2068661122070742
$ws.Range("J8").Value = 0.2068661122070743
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 27.73243066666667
$ws.Range("N8").Value = 83.197292
$ws.Range("O8").Value = 0.3628105447549136
$ws.Range("P8").Value = 0.3628105447549136
$ws.Range("Q8").Value = 8.498427739072445
$ws.Range("R8").Value = 76.485849651652
$ws.Range("S8").Value = 0.07505320686117968
$ws.Range("T8").Value = 0.0750532068611797

$ws.Range("E9").Value = 2
$ws.Range("F9").Value = 0.6666666666666666
$ws.Range("G9").Value = 0.3064436666666667
$ws.Range("H9").Value = 0.919331
$ws.Range("I9").Value = 0.2068661122070742
$ws.Range("J9").Value = 0.2068661122070743
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 1.096171666666667
$ws.Range("N9").Value = 3.288515
$ws.Range("O9").Value = 0.01434070616847367
$ws.Range("P9").Value = 0.01434070616847367
$ws.Range("Q9").Value = 0.3359148648294445
$ws.Range("R9").Value = 3.023233783465
$ws.Range("S9").Value = 0.002966606131376155
$ws.Range("T9").Value = 0.002966606131376155

$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 0.5451493333333333
$ws.Range("H10").Value = 1.635448
$ws.Range("I10").Value = 0.3680053968340403
$ws.Range("J10").Value = 0.3680053968340404
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 6.322177333333333
$ws.Range("N10").Value = 18.966532
$ws.Range("O10").Value = 0.08271011762055308
$ws.Range("P10").Value = 0.0827101176205531
$ws.Range("Q10").Value = 3.446530758481777
$ws.Range("R10").Value = 31.018776826336
$ws.Range("S10").Value = 0.03043776965714179
$ws.Range("T10").Value = 0.0304377696571418

$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 0.5451493333333333
$ws.Range("H11").Value = 1.635448
$ws.Range("I11").Value = 0.3680053968340403
$ws.Range("J11").Value = 0.3680053968340404
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 41.286995
$ws.Range("N11").Value = 123.860985
$ws.Range("O11").Value = 0.5401386314560596
$ws.Range("P11").Value = 0.5401386314560597
$ws.Range("Q11").Value = 22.50757779958666
$ws.Range("R11").Value = 202.56820019628
$ws.Range("S11").Value = 0.1987739314143827
$ws.Range("T11").Value = 0.1987739314143827

$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 0.5451493333333333
$ws.Range("H12").Value = 1.635448
$ws.Range("I12").Value = 0.3680053968340403
$ws.Range("J12").Value = 0.3680053968340404
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 27.73243066666667
$ws.Range("N12").Value = 83.197292
$ws.Range("O12").Value = 0.3628105447549136
$ws.Range("P12").Value = 0.3628105447549136
$ws.Range("Q12").Value = 15.11831608964622
$ws.Range("R12").Value = 136.064844806816
$ws.Range("S12").Value = 0.1335162384981063
$ws.Range("T12").Value = 0.1335162384981063

$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 0.5451493333333333
$ws.Range("H13").Value = 1.635448
$ws.Range("I13").Value = 0.3680053968340403
$ws.Range("J13").Value = 0.3680053968340404
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 1.096171666666667
$ws.Range("N13").Value = 3.288515
$ws.Range("O13").Value = 0.01434070616847367
$ws.Range("P13").Value = 0.01434070616847367
$ws.Range("Q13").Value = 0.5975772533022222
$ws.Range("R13").Value = 5.37819527972
$ws.Range("S13").Value = 0.005277457264409521
$ws.Range("T13").Value = 0.005277457264409522

$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.05723766666666667
$ws.Range("H14").Value = 0.171713
$ws.Range("I14").Value = 0.03863853250397663
$ws.Range("J14").Value = 0.03863853250397663
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 6.322177333333333
$ws.Range("N14").Value = 18.966532
$ws.Range("O14").Value = 0.08271011762055308
$ws.Range("P14").Value = 0.0827101176205531
$ws.Range("Q14").Value = 0.3618666788128889
$ws.Range("R14").Value = 3.256800109316
$ws.Range("S14").Value = 0.00319579756808947
$ws.Range("T14").Value = 0.003195797568089472

$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.05723766666666667
$ws.Range("H15").Value = 0.171713
$ws.Range("I15").Value = 0.03863853250397663
$ws.Range("J15").Value = 0.03863853250397663
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 41.286995
$ws.Range("N15").Value = 123.860985
$ws.Range("O15").Value = 0.5401386314560596
$ws.Range("P15").Value = 0.5401386314560597
$ws.Range("Q15").Value = 2.363171257478333
$ws.Range("R15").Value = 21.268541317305
$ws.Range("S15").Value = 0.02087016406816841
$ws.Range("T15").Value = 0.02087016406816842

$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.05723766666666667
$ws.Range("H16").Value = 0.171713
$ws.Range("I16").Value = 0.03863853250397663
$ws.Range("J16").Value = 0.03863853250397663
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 27.73243066666667
$ws.Range("N16").Value = 83.197292
$ws.Range("O16").Value = 0.3628105447549136
$ws.Range("P16").Value = 0.3628105447549136
$ws.Range("Q16").Value = 1.587339622355111
$ws.Range("R16").Value = 14.286056601196
$ws.Range("S16").Value = 0.0140184670262982
$ws.Range("T16").Value = 0.0140184670262982

$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.05723766666666667
$ws.Range("H17").Value = 0.171713
$ws.Range("I17").Value = 0.03863853250397663
$ws.Range("J17").Value = 0.03863853250397663
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 1.096171666666667
$ws.Range("N17").Value = 3.288515
$ws.Range("O17").Value = 0.01434070616847367
$ws.Range("P17").Value = 0.01434070616847367
$ws.Range("Q17").Value = 0.06274230846611112
$ws.Range("R17").Value = 0.564680776195
$ws.Range("S17").Value = 0.0005541038414205478
$ws.Range("T17").Value = 0.000554103841420548

